$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demi")

# --- Shift existing rows 18-20 down to 19-21 (copy, bottom-up, to preserve values+formatting) ---
$ws.Range("A20:D20").Copy($ws.Range("A21:D21"))
$ws.Range("A19:D19").Copy($ws.Range("A20:D20"))
$ws.Range("A18:D18").Copy($ws.Range("A19:D19"))

# --- Prepare new row 22 (copy formatting for A/B/C from row 21, which now holds the former row 20's
#     formatting - this keeps B22's date number format identical to the other date cells, and avoids
#     creating a brand-new style entry) ---
$ws.Range("A21").Copy($ws.Range("A22"))
$ws.Range("B21").Copy($ws.Range("B22"))
$ws.Range("C21").Copy($ws.Range("C22"))
$ws.Range("A22:C22").ClearContents()

# --- New row 18 content: "Thuis alleen gewerkt" on 27/11/20203 (typo'd text date), 60 minutes ---
$ws.Range("A18:D18").ClearContents()
$ws.Range("B18").ClearFormats()
$ws.Range("A18").Value2 = "Thuis alleen gewerkt"
$ws.Range("B18").Value2 = "27/11/20203"
$ws.Range("C18").Value2 = 60

# --- New row 22 content: "Thuis alleen gewerkt" on 12/1/2023, 120 minutes ---
$ws.Range("A22").Value2 = "Thuis alleen gewerkt"
$ws.Range("B22").Value2 = 45261
$ws.Range("C22").Value2 = 120

# --- Make "Demi" the active sheet/tab, scrolled + zoomed to where the new rows were added, with the
#     last-entered cell (D22) selected ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 130
$win.ScrollRow = 8
$win.ScrollColumn = 1
$ws.Range("D22").Select()
